$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column Q mirrors column P's formatting (2020 data column), row by row.
$ws.Range("P3").Copy()
$ws.Range("Q3").PasteSpecial(-4122)

$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)
$ws.Range("Q4").Value = 2020

$ws.Range("P5").Copy()
$ws.Range("Q5").PasteSpecial(-4122)
$ws.Range("Q5").Value = 38.6

$ws.Range("P6").Copy()
$ws.Range("Q6").PasteSpecial(-4122)
$ws.Range("Q6").Value = 42.4

$ws.Range("P7").Copy()
$ws.Range("Q7").PasteSpecial(-4122)
$ws.Range("Q7").Value = 53.2

$ws.Range("P8").Copy()
$ws.Range("Q8").PasteSpecial(-4122)
$ws.Range("Q8").Value = 90.6

$ws.Range("P9").Copy()
$ws.Range("Q9").PasteSpecial(-4122)
$ws.Range("Q9").Value = 52.6

$ws.Range("P10").Copy()
$ws.Range("Q10").PasteSpecial(-4122)
$ws.Range("Q10").Value = 24.5

$ws.Range("P11").Copy()
$ws.Range("Q11").PasteSpecial(-4122)
$ws.Range("Q11").Value = 69.1

$ws.Range("P12").Copy()
$ws.Range("Q12").PasteSpecial(-4122)
$ws.Range("Q12").Value = 32.2

$ws.Range("P13").Copy()
$ws.Range("Q13").PasteSpecial(-4122)
$ws.Range("Q13").Value = 19.1

$ws.Range("P14").Copy()
$ws.Range("Q14").PasteSpecial(-4122)
$ws.Range("Q14").Value = 25.2

$ws.Range("R27").Select()
